# Bondtech Kit reference added
# Rows 3-8 (all BMG/Bondtech "BUY" parts) now reference the combined
# "BMG Internals Set for HextrudORT" kit instead of individual Bondtech
# product pages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$comment = "Included in BMG Internals Set for HextrudORT"
$makeBuy = "(BUY) KIT"
$kitUrl  = "https://www.bondtech.se/product/bmg-internals-set-for-hextrudort/"

$rows = @(3, 4, 5, 6, 7, 8)

foreach ($r in $rows) {
    $ws.Range("F$r").Value = $comment
    $ws.Range("G$r").Value = $makeBuy
    $ws.Range("K$r").Value = $kitUrl
}

# Qty for the bearing line (row 8) doubles since it now comes from the kit
$ws.Range("H8").Value = "2*"

# Replace the individual Bondtech product hyperlinks with the single kit
# hyperlink (each K3..K8 cell keeps its own hyperlink object, just now
# pointing at the combined kit page).
foreach ($hl in $ws.Hyperlinks) {
    $hl.Address = $kitUrl
}

# Selection moved from O6 to I4 in the saved file
$ws.Range("I4").Select()
